$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7 = 235
    8 = 57
    9 = 258
    10 = 269
    11 = 728
    12 = 940
    13 = 558
    14 = 740
    15 = 629
    16 = 367
    17 = 391
    18 = 175
    19 = 12
    20 = 20
    21 = 154
    22 = 360
    23 = 156
    24 = 155
    25 = 458
    26 = 467
    27 = 98
    28 = 62
    29 = 700
    30 = 515
    31 = 697
    32 = 602
    33 = 105
    34 = 102
    35 = 338
    36 = 207
    37 = 548
    38 = 343
    39 = 80
    40 = 172
    41 = 402
    42 = 342
    43 = 961
    44 = 1108
    45 = 197
    46 = 161
    47 = 1750
    48 = 1794
    49 = 314
    50 = 274
    51 = 696
    52 = 742
    53 = 470
    54 = 308
    55 = 728
    56 = 597
    57 = 1347
    58 = 1605
}

foreach ($row in $updates.Keys) {
    $ws.Range("O$row").Value = $updates[$row]
}
